$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Config2 worksheet: fill in the "Hawkeye" (row 25) and "OPTGen" (row 26)
# results for the GemsFDTD benchmark, which previously only had the
# benchmark/policy labels (columns A/B) with #DIV/0! placeholders in H/I.
# (Done first so the workbook's active sheet ends up back on Config1, as
# in the original file.)
# ---------------------------------------------------------------------------
$wsConfig2 = $wb.Worksheets.Item("Config2")

# Row 25 - Hawkeye
$wsConfig2.Range("C25").Value = 50000001
$wsConfig2.Range("D25").Value = 80542122
$wsConfig2.Range("E25").Value = 2703163
$wsConfig2.Range("F25").Value = 101472
$wsConfig2.Range("G25").Value = 2601691

# Row 26 - OPTGen
$wsConfig2.Range("C26").Value = 50000001
$wsConfig2.Range("D26").Value = 80542122
$wsConfig2.Range("E26").Value = 9064
$wsConfig2.Range("F26").Value = 2963
$wsConfig2.Range("G26").Formula = "=E26-F26"
$wsConfig2.Range("J26").Formula = "=F26/E26"

# Update the stored selection for the Config2 sheet.
$wsConfig2.Activate()
$wsConfig2.Range("C26").Select()

# ---------------------------------------------------------------------------
# Config1 worksheet: same fill-in for its own GemsFDTD Hawkeye / OPTGen rows.
# ---------------------------------------------------------------------------
$wsConfig1 = $wb.Worksheets.Item("Config1")

# Row 25 - Hawkeye
$wsConfig1.Range("C25").Value = 50000001
$wsConfig1.Range("D25").Value = 151789738
$wsConfig1.Range("E25").Value = 2633225
$wsConfig1.Range("F25").Value = 582227
$wsConfig1.Range("G25").Value = 2050998

# Row 26 - OPTGen
$wsConfig1.Range("C26").Value = 50000001
$wsConfig1.Range("D26").Value = 151789738
$wsConfig1.Range("E26").Value = 47832
$wsConfig1.Range("F26").Value = 12892
$wsConfig1.Range("G26").Formula = "=E26-F26"
$wsConfig1.Range("J26").Formula = "=F26/E26"

# Update the stored selection for the Config1 sheet and make it active again,
# matching the original workbook's active tab.
$wsConfig1.Activate()
$wsConfig1.Range("C28").Select()
